# Insert a new data row at row 348 (pushes the existing row 348..431 down to
# 349..432, automatically growing the sheet dimension from A1:R431 to
# A1:R432), then populate the newly inserted row with the new "Repollo"
# price entry for Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 348-431 down to 349-432 by inserting a blank row at 348.
$ws.Rows.Item(348).Insert()

# Populate the new row 348 with the new record's data.
$ws.Range("A348").Value = 3
$ws.Range("B348").Value = "Femacal de La Calera"
$ws.Range("C348").Value = "Coquimbo"
$ws.Range("D348").Value = 44543
$ws.Range("E348").Value = 5
$ws.Range("F348").Value = 100112006
$ws.Range("G348").Value = "Repollo"
$ws.Range("H348").Value = "Crespo record"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 1600
$ws.Range("K348").Value = 600
$ws.Range("L348").Value = 600
$ws.Range("M348").Value = 600
$ws.Range("N348").Value = '$/unidad'
$ws.Range("O348").Value = "Provincia de Quillota"
$ws.Range("P348").Value = 600
$ws.Range("Q348").Value = 1
$ws.Range("R348").Value = "Hortaliza"
